$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-03-13 Wednesday" "2024-03-14 Thursday"
Replace-Text "15÷9=1, 6" "31÷7=4, 3"
Replace-Text "72÷4=18, 0" "48÷3=16, 0"
Replace-Text "74÷9=8, 2" "64÷9=7, 1"
Replace-Text "83÷8=10, 3" "18÷6=3, 0"
Replace-Text "77÷5=15, 2" "60÷3=20, 0"
Replace-Text "40÷6=6, 4" "44÷5=8, 4"
Replace-Text "21÷8=2, 5" "19÷8=2, 3"
Replace-Text "68÷6=11, 2" "75÷9=8, 3"
Replace-Text "46÷9=5, 1" "39÷4=9, 3"
Replace-Text "91÷7=13, 0" "95÷7=13, 4"
Replace-Text "13÷9=1, 4" "34÷9=3, 7"
Replace-Text "99÷2=49, 1" "73÷5=14, 3"
Replace-Text "67÷9=7, 4" "25÷8=3, 1"
Replace-Text "87÷2=43, 1" "62÷2=31, 0"
Replace-Text "83÷6=13, 5" "17÷4=4, 1"
Replace-Text "21÷6=3, 3" "95÷2=47, 1"
Replace-Text "82÷5=16, 2" "64÷5=12, 4"
Replace-Text "64÷7=9, 1" "79÷4=19, 3"
Replace-Text "25÷7=3, 4" "41÷6=6, 5"
Replace-Text "28÷2=14, 0" "70÷3=23, 1"
Replace-Text "74÷4=18, 2" "74÷8=9, 2"
Replace-Text "18÷2=9, 0" "50÷7=7, 1"
Replace-Text "87÷7=12, 3" "25÷8=3, 1"
Replace-Text "36÷3=12, 0" "35÷5=7, 0"
Replace-Text "52÷8=6, 4" "84÷3=28, 0"
